# ZBP_03_strategie_domacnosti.xlsx - add a new survey wave (30. 3. 2021)
# as the next column on both sheets, and bump the "aktualizace" date in
# the two footer/title strings from 23. 3. 2021 to 7. 4. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": percentages table. New column AA holds the 30. 3. 2021
# wave, added right after column Z (16. 3. 2021).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Give the new header cell the same look (bold, centered, bordered) as the
# rest of row 1 before filling in its text.
$wsData.Range("Z1").Copy()
$wsData.Range("AA1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsData.Range("AA1").Value = "30. 3. 2021"

$dataValues = @{
    2  = 0.22
    3  = 0.12
    4  = 0.54
    5  = 0.28
    6  = 0.14
    7  = 0.21
    8  = 0.26
    9  = 0.2
    10 = 0.26
    11 = 0.21
    12 = 0.21
    13 = 0.38
    14 = 0.21
    15 = 0.21
    16 = 0.25
    17 = 0.19
    18 = 0.28
    19 = 0.29
    20 = 0.16
    21 = 0.14
    22 = 0.13
    23 = 0.21
    24 = 0.45
    25 = 0.41
    26 = 0.11
    27 = 0.07000000000000001
    28 = 0.12
    29 = 0.18
    30 = 0.09
    31 = 0.11
    32 = 0.12
    33 = 0.16
    34 = 0.19
    35 = 0.1
    36 = 0.13
    37 = 0.13
    38 = 0.07000000000000001
    39 = 0.26
    40 = 0.14
    41 = 0.07000000000000001
    42 = 0.05
    43 = 0.05
    44 = 0.15
    45 = 0.3
}

foreach ($row in $dataValues.Keys) {
    $wsData.Cells.Item($row, 27).Value = $dataValues[$row]
}

# Footer title string: bump the "aktualizace" date.
$wsData.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 7. 4. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": respondent-count table. New column Z holds the
# 30. 3. 2021 wave, added right after column Y (16. 3. 2021).
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# Same header styling treatment as on the "data" sheet.
$wsPocet.Range("Y1").Copy()
$wsPocet.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsPocet.Range("Z1").Value = "30. 3. 2021"

$pocetValues = @{
    2  = 2061
    3  = 228
    4  = 458
    5  = 1375
    6  = 983
    7  = 183
    8  = 590
    9  = 305
    10 = 936
    11 = 170
    12 = 134
    13 = 821
    14 = 942
    15 = 712
    16 = 407
    17 = 251
    18 = 761
    19 = 649
    20 = 254
    21 = 551
    22 = 360
    23 = 231
}

foreach ($row in $pocetValues.Keys) {
    $wsPocet.Cells.Item($row, 26).Value = $pocetValues[$row]
}

# Footer title string: bump the "aktualizace" date.
$wsPocet.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 7. 4. 2021"

# The footer row's other cells (B24:Y24) are blank placeholders; keep the
# newly added Z24 blank too, matching the row's existing pattern.
$wsPocet.Range("B24:Z24").Value = ""
